# V 1.0.4: Fix AGIP
#
# The "Clientes" sheet's T column (Resultado AGIP) looks up each row's
# number in the external "Control.xlsx" workbook's AGIP sheet via
# IFERROR(VLOOKUP(...)). The external workbook's AGIP cache was cleared so
# those lookups no longer resolve (mirroring the ARBA/"S" column, which
# already comes back blank). Re-touching/recalculating the sheet is what
# flips the stale cached "OK" / "No Hay Retenciones para el Período"
# results back to blank, same as the real fix.
#
# The user had also scrolled/re-selected the sheet before saving (view was
# at H2, ended up at H13 with column D pinned to the left edge) - reproduce
# that as closely as the object model allows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")
$ws.Activate()

# Force the external-reference formulas in column T (and S) to recompute
# against the refreshed external link cache, clearing out the stale
# cached results.
try { $wb.RefreshAll() } catch {}
try { $wb.UpdateLink("Control.xlsx", 5) } catch {}
$excel.CalculateFull()

$win = $excel.ActiveWindow

# Scroll so column D becomes the left-most visible column (best effort -
# matches the saved view's topLeftCell="D1").
try { $win.ScrollColumn = 4 } catch {}
try { $win.ScrollRow = 1 } catch {}

# Leave the selection where the file was saved: H13.
$ws.Range("H13").Select()
